# Update "想去人数" (F column) figures (and one "最低票价" G column that
# flipped from a numeric price to "已售罄" / Sold Out) across the three
# affected sheets: 展览, 演出, 全部类型. (本地生活 is untouched.)

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 2859
$ws.Range("F7").Value = 1686
$ws.Range("F8").Value = 1871
$ws.Range("F11").Value = 763
$ws.Range("F12").Value = 912
$ws.Range("F14").Value = 382
$ws.Range("F17").Value = 52
$ws.Range("F19").Value = 6796
$ws.Range("F21").Value = 1617
$ws.Range("F22").Value = 166
$ws.Range("F23").Value = 181
$ws.Range("F25").Value = 313
$ws.Range("F26").Value = 273
$ws.Range("F27").Value = 72
$ws.Range("F28").Value = 1108
$ws.Range("F31").Value = 96
$ws.Range("F33").Value = 791
$ws.Range("F34").Value = 1929
$ws.Range("G34").Value = "已售罄"
$ws.Range("F35").Value = 164
$ws.Range("F36").Value = 146
$ws.Range("F37").Value = 230
$ws.Range("F38").Value = 26
$ws.Range("F40").Value = 228
$ws.Range("F41").Value = 78
$ws.Range("F42").Value = 168

# ---- Sheet: 演出 (Performance) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 4
$ws.Range("F9").Value = 8

# ---- Sheet: 全部类型 (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 2859
$ws.Range("F10").Value = 1686
$ws.Range("F11").Value = 1871
$ws.Range("F14").Value = 763
$ws.Range("F16").Value = 912
$ws.Range("F18").Value = 382
$ws.Range("F20").Value = 52
$ws.Range("F22").Value = 6796
$ws.Range("F24").Value = 1617
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 166
$ws.Range("F27").Value = 181
$ws.Range("F29").Value = 313
$ws.Range("F30").Value = 273
$ws.Range("F31").Value = 72
$ws.Range("F32").Value = 1108
$ws.Range("F35").Value = 96
$ws.Range("F37").Value = 791
$ws.Range("F38").Value = 1929
$ws.Range("G38").Value = "已售罄"
$ws.Range("F39").Value = 164
$ws.Range("F40").Value = 146
$ws.Range("F41").Value = 230
$ws.Range("F42").Value = 26
$ws.Range("F44").Value = 228
$ws.Range("F47").Value = 8
$ws.Range("F48").Value = 78
$ws.Range("F49").Value = 168
